$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are stored as literal text in the source data
# (even when they look numeric), so force text format before assigning,
# then clear the format override to keep the original (default) cell style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.847.67'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.02'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.48'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5346'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3751'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07183'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.59'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08150'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.891.81'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.17'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.317'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.84'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008536'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.894.89'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.981'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.71%  '
$ws.Range("E22").Value = '  -1.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.399'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.16'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.09'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.01'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.726'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.607'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09143'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8124'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05013'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.174'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.949'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6018'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.214'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.609'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.21%  '
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.069'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.629'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.911'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.06'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5095'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1492'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.22%  '
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.942'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.08%  '
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.69'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06053'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("E51").Value = '  -2.84%  '
